$wb = $excel.ActiveWorkbook

$changes = @(
    @{Sheet="Citywide Totals"; Cell="L2"; Value=6444},
    @{Sheet="Citywide Totals"; Cell="L3"; Value=6950},
    @{Sheet="Citywide Totals"; Cell="J4"; Value=1884},
    @{Sheet="Citywide Totals"; Cell="K4"; Value=1791},
    @{Sheet="Citywide Totals"; Cell="L4"; Value=1732},
    @{Sheet="Citywide Totals"; Cell="L5"; Value=408},
    @{Sheet="Citywide Totals"; Cell="L6"; Value=5700},
    @{Sheet="Citywide Totals"; Cell="J7"; Value=29362},
    @{Sheet="Citywide Totals"; Cell="K7"; Value=27585},
    @{Sheet="Citywide Totals"; Cell="L7"; Value=21234},
    @{Sheet="Logan Square"; Cell="L2"; Value=72},
    @{Sheet="Logan Square"; Cell="K4"; Value=19},
    @{Sheet="Logan Square"; Cell="L6"; Value=78},
    @{Sheet="Logan Square"; Cell="K7"; Value=350},
    @{Sheet="Logan Square"; Cell="L7"; Value=240},
    @{Sheet="Austin"; Cell="L2"; Value=430},
    @{Sheet="Austin"; Cell="L3"; Value=494},
    @{Sheet="Austin"; Cell="L6"; Value=339},
    @{Sheet="Austin"; Cell="L7"; Value=1404},
    @{Sheet="South Chicago"; Cell="L2"; Value=149},
    @{Sheet="South Chicago"; Cell="L4"; Value=18},
    @{Sheet="South Chicago"; Cell="L7"; Value=467},
    @{Sheet="Garfield Park"; Cell="L3"; Value=338},
    @{Sheet="Garfield Park"; Cell="L5"; Value=24},
    @{Sheet="Garfield Park"; Cell="L7"; Value=957},
    @{Sheet="Grand Crossing"; Cell="L2"; Value=248},
    @{Sheet="Grand Crossing"; Cell="L3"; Value=289},
    @{Sheet="Grand Crossing"; Cell="L6"; Value=210},
    @{Sheet="Grand Crossing"; Cell="L7"; Value=816},
    @{Sheet="New City"; Cell="L2"; Value=153},
    @{Sheet="New City"; Cell="L3"; Value=132},
    @{Sheet="New City"; Cell="L7"; Value=417},
    @{Sheet="Fuller Park"; Cell="L2"; Value=26},
    @{Sheet="Fuller Park"; Cell="L7"; Value=94},
    @{Sheet="By Neighborhood"; Cell="L2"; Value=187},
    @{Sheet="By Neighborhood"; Cell="L5"; Value=76},
    @{Sheet="By Neighborhood"; Cell="L6"; Value=171},
    @{Sheet="By Neighborhood"; Cell="L7"; Value=673},
    @{Sheet="By Neighborhood"; Cell="L8"; Value=1404},
    @{Sheet="By Neighborhood"; Cell="L15"; Value=178},
    @{Sheet="By Neighborhood"; Cell="L19"; Value=582},
    @{Sheet="By Neighborhood"; Cell="L20"; Value=540},
    @{Sheet="By Neighborhood"; Cell="L21"; Value=69},
    @{Sheet="By Neighborhood"; Cell="L27"; Value=182},
    @{Sheet="By Neighborhood"; Cell="L29"; Value=1190},
    @{Sheet="By Neighborhood"; Cell="L30"; Value=94},
    @{Sheet="By Neighborhood"; Cell="L31"; Value=212},
    @{Sheet="By Neighborhood"; Cell="L33"; Value=957},
    @{Sheet="By Neighborhood"; Cell="L35"; Value=27},
    @{Sheet="By Neighborhood"; Cell="K36"; Value=357},
    @{Sheet="By Neighborhood"; Cell="L36"; Value=271},
    @{Sheet="By Neighborhood"; Cell="L37"; Value=816},
    @{Sheet="By Neighborhood"; Cell="L42"; Value=668},
    @{Sheet="By Neighborhood"; Cell="L48"; Value=277},
    @{Sheet="By Neighborhood"; Cell="L51"; Value=263},
    @{Sheet="By Neighborhood"; Cell="K53"; Value=350},
    @{Sheet="By Neighborhood"; Cell="L53"; Value=240},
    @{Sheet="By Neighborhood"; Cell="L59"; Value=36},
    @{Sheet="By Neighborhood"; Cell="J63"; Value=238},
    @{Sheet="By Neighborhood"; Cell="K63"; Value=179},
    @{Sheet="By Neighborhood"; Cell="L63"; Value=71},
    @{Sheet="By Neighborhood"; Cell="L64"; Value=135},
    @{Sheet="By Neighborhood"; Cell="L65"; Value=417},
    @{Sheet="By Neighborhood"; Cell="L66"; Value=61},
    @{Sheet="By Neighborhood"; Cell="L67"; Value=734},
    @{Sheet="By Neighborhood"; Cell="L71"; Value=53},
    @{Sheet="By Neighborhood"; Cell="L72"; Value=84},
    @{Sheet="By Neighborhood"; Cell="L73"; Value=167},
    @{Sheet="By Neighborhood"; Cell="L75"; Value=77},
    @{Sheet="By Neighborhood"; Cell="L76"; Value=334},
    @{Sheet="By Neighborhood"; Cell="L78"; Value=281},
    @{Sheet="By Neighborhood"; Cell="L79"; Value=588},
    @{Sheet="By Neighborhood"; Cell="L83"; Value=467},
    @{Sheet="By Neighborhood"; Cell="L85"; Value=1051},
    @{Sheet="By Neighborhood"; Cell="L86"; Value=133},
    @{Sheet="By Neighborhood"; Cell="L88"; Value=225},
    @{Sheet="By Neighborhood"; Cell="L89"; Value=286},
    @{Sheet="By Neighborhood"; Cell="L91"; Value=286},
    @{Sheet="By Neighborhood"; Cell="L94"; Value=258},
    @{Sheet="By Neighborhood"; Cell="L96"; Value=233},
    @{Sheet="By Neighborhood"; Cell="J101"; Value=29362},
    @{Sheet="By Neighborhood"; Cell="K101"; Value=27585},
    @{Sheet="By Neighborhood"; Cell="L101"; Value=21234},
    @{Sheet="Gage Park"; Cell="L4"; Value=13},
    @{Sheet="Gage Park"; Cell="L6"; Value=55},
    @{Sheet="Gage Park"; Cell="L7"; Value=212},
    @{Sheet="North Lawndale"; Cell="L3"; Value=286},
    @{Sheet="North Lawndale"; Cell="L7"; Value=734},
    @{Sheet="Englewood"; Cell="L2"; Value=363},
    @{Sheet="Englewood"; Cell="L7"; Value=1190},
    @{Sheet="Lake View"; Cell="L2"; Value=42},
    @{Sheet="Lake View"; Cell="L7"; Value=277},
    @{Sheet="Chatham"; Cell="L4"; Value=29},
    @{Sheet="Chatham"; Cell="L6"; Value=160},
    @{Sheet="Chatham"; Cell="L7"; Value=582},
    @{Sheet="River North"; Cell="L6"; Value=151},
    @{Sheet="River North"; Cell="L7"; Value=334},
    @{Sheet="Ashburn"; Cell="L6"; Value=34},
    @{Sheet="Ashburn"; Cell="L7"; Value=171},
    @{Sheet="Humboldt Park"; Cell="L2"; Value=182},
    @{Sheet="Humboldt Park"; Cell="L7"; Value=668},
    @{Sheet="Rogers Park"; Cell="L2"; Value=73},
    @{Sheet="Rogers Park"; Cell="L6"; Value=81},
    @{Sheet="Rogers Park"; Cell="L7"; Value=281},
    @{Sheet="West Ridge"; Cell="L6"; Value=66},
    @{Sheet="West Ridge"; Cell="L7"; Value=233},
    @{Sheet="Washington Park"; Cell="L3"; Value=128},
    @{Sheet="Washington Park"; Cell="L4"; Value=16},
    @{Sheet="Washington Park"; Cell="L6"; Value=37},
    @{Sheet="Washington Park"; Cell="L7"; Value=286},
    @{Sheet="Chinatown"; Cell="L3"; Value=19},
    @{Sheet="Chinatown"; Cell="L7"; Value=69},
    @{Sheet="Roseland"; Cell="L6"; Value=157},
    @{Sheet="Roseland"; Cell="L7"; Value=588},
    @{Sheet="Near South Side"; Cell="L2"; Value=39},
    @{Sheet="Near South Side"; Cell="L7"; Value=135},
    @{Sheet="Chicago Lawn"; Cell="L2"; Value=170},
    @{Sheet="Chicago Lawn"; Cell="L3"; Value=189},
    @{Sheet="Chicago Lawn"; Cell="L7"; Value=540},
    @{Sheet="Grand Boulevard"; Cell="L2"; Value=92},
    @{Sheet="Grand Boulevard"; Cell="L3"; Value=89},
    @{Sheet="Grand Boulevard"; Cell="K5"; Value=7},
    @{Sheet="Grand Boulevard"; Cell="L6"; Value=66},
    @{Sheet="Grand Boulevard"; Cell="K7"; Value=357},
    @{Sheet="Grand Boulevard"; Cell="L7"; Value=271},
    @{Sheet="Auburn Gresham"; Cell="L3"; Value=216},
    @{Sheet="Auburn Gresham"; Cell="L7"; Value=673},
    @{Sheet="West Loop"; Cell="L3"; Value=62},
    @{Sheet="West Loop"; Cell="L7"; Value=258},
    @{Sheet="Brighton Park"; Cell="L4"; Value=16},
    @{Sheet="Brighton Park"; Cell="L7"; Value=178},
    @{Sheet="North Center"; Cell="L3"; Value=19},
    @{Sheet="North Center"; Cell="L7"; Value=61},
    @{Sheet="Gold Coast"; Cell="L6"; Value=12},
    @{Sheet="Gold Coast"; Cell="L7"; Value=27},
    @{Sheet="Portage Park"; Cell="L6"; Value=38},
    @{Sheet="Portage Park"; Cell="L7"; Value=167},
    @{Sheet="Montclare"; Cell="L3"; Value=16},
    @{Sheet="Montclare"; Cell="L7"; Value=36},
    @{Sheet="Albany Park"; Cell="L6"; Value=48},
    @{Sheet="Albany Park"; Cell="L7"; Value=187},
    @{Sheet="United Center"; Cell="L3"; Value=76},
    @{Sheet="United Center"; Cell="L6"; Value=62},
    @{Sheet="United Center"; Cell="L7"; Value=225},
    @{Sheet="Uptown"; Cell="L3"; Value=85},
    @{Sheet="Uptown"; Cell="L7"; Value=286},
    @{Sheet="Armour Square"; Cell="L2"; Value=17},
    @{Sheet="Armour Square"; Cell="L7"; Value=76},
    @{Sheet="Edgewater"; Cell="L6"; Value=57},
    @{Sheet="Edgewater"; Cell="L7"; Value=182},
    @{Sheet="Streeterville"; Cell="L6"; Value=17},
    @{Sheet="Streeterville"; Cell="L7"; Value=133},
    @{Sheet="Pullman"; Cell="L3"; Value=27},
    @{Sheet="Pullman"; Cell="L7"; Value=77},
    @{Sheet="Little Italy, UIC"; Cell="L4"; Value=39},
    @{Sheet="Little Italy, UIC"; Cell="L7"; Value=263},
    @{Sheet="South Shore"; Cell="L3"; Value=438},
    @{Sheet="South Shore"; Cell="L6"; Value=216},
    @{Sheet="South Shore"; Cell="L7"; Value=1051},
    @{Sheet="Oakland"; Cell="L3"; Value=17},
    @{Sheet="Oakland"; Cell="L7"; Value=53},
    @{Sheet="Old Town"; Cell="L6"; Value=27},
    @{Sheet="Old Town"; Cell="L7"; Value=84}
)

foreach ($ch in $changes) {
    $ws = $wb.Worksheets.Item($ch.Sheet)
    $ws.Range($ch.Cell).Value = $ch.Value
}

Write-Output "Applied $($changes.Count) cell updates"
